$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fig32")

# Update title and source strings: January 2017 -> February 2017
$ws.Range("A2").Value = "Short-Term Energy Outlook, February 2017"
$ws.Range("A56").Value = "Source: Short-Term Energy Outlook, February 2017."

# Update the World production (C), World consumption (D), and implied
# stock change and balance (E) data series for rows 28-55 (2012-Q1..2018-Q4)
# Values are cast via [double] on the literal string to avoid PowerShell
# parser issues with scientific-notation literals (e.g. "E-3").
$ws.Range("C28").Value = [double]"90.615488103999994"
$ws.Range("D28").Value = [double]"89.876435060999995"
$ws.Range("E28").Value = [double]"0.73905304256000004"

$ws.Range("C29").Value = [double]"90.450962777000001"
$ws.Range("D29").Value = [double]"90.284702682000002"
$ws.Range("E29").Value = [double]"0.16626009423999999"

$ws.Range("C30").Value = [double]"90.557449757000001"
$ws.Range("D30").Value = [double]"91.632804859999993"
$ws.Range("E30").Value = [double]"-1.0753551027999999"

$ws.Range("C31").Value = [double]"91.001595703999996"
$ws.Range("D31").Value = [double]"92.568502902000006"
$ws.Range("E31").Value = [double]"-1.5669071974"

$ws.Range("C32").Value = [double]"90.049844258999997"
$ws.Range("D32").Value = [double]"91.167427353999997"
$ws.Range("E32").Value = [double]"-1.1175830953999999"

$ws.Range("C33").Value = [double]"91.249070962000005"
$ws.Range("D33").Value = [double]"91.851533747000005"
$ws.Range("E33").Value = [double]"-0.60246278483000004"

$ws.Range("C34").Value = [double]"91.814030509000006"
$ws.Range("D34").Value = [double]"93.036322092000006"
$ws.Range("E34").Value = [double]"-1.2222915831000001"

$ws.Range("C35").Value = [double]"91.863069675999995"
$ws.Range("D35").Value = [double]"93.217316842000002"
$ws.Range("E35").Value = [double]"-1.3542471661"

$ws.Range("C36").Value = [double]"92.256456635000006"
$ws.Range("D36").Value = [double]"92.872456141000001"
$ws.Range("E36").Value = [double]"-0.61599950595999997"

$ws.Range("C37").Value = [double]"92.940610668999994"
$ws.Range("D37").Value = [double]"92.945229963000003"
$ws.Range("E37").Value = [double]"-4.6192939900000004E-3"

$ws.Range("C38").Value = [double]"94.310605315999993"
$ws.Range("D38").Value = [double]"93.985852453000007"
$ws.Range("E38").Value = [double]"0.32475286288999999"

$ws.Range("C39").Value = [double]"95.780064386999996"
$ws.Range("D39").Value = [double]"94.590933238000005"
$ws.Range("E39").Value = [double]"1.1891311483"

$ws.Range("C40").Value = [double]"95.573163984000004"
$ws.Range("D40").Value = [double]"94.057179778999995"
$ws.Range("E40").Value = [double]"1.5159842051000001"

$ws.Range("C41").Value = [double]"96.559370178999998"
$ws.Range("D41").Value = [double]"94.595104372999998"
$ws.Range("E41").Value = [double]"1.9642658061"

$ws.Range("C42").Value = [double]"97.485330544999997"
$ws.Range("D42").Value = [double]"96.034932707999999"
$ws.Range("E42").Value = [double]"1.4503978367999999"

$ws.Range("C43").Value = [double]"97.670774522000002"
$ws.Range("D43").Value = [double]"95.526807508000005"
$ws.Range("E43").Value = [double]"2.1439670141999998"

$ws.Range("C44").Value = [double]"96.823780651999996"
$ws.Range("D44").Value = [double]"95.360327831999996"
$ws.Range("E44").Value = [double]"1.4634528198000001"

$ws.Range("C45").Value = [double]"96.489593322999994"
$ws.Range("D45").Value = [double]"96.087232920999995"
$ws.Range("E45").Value = [double]"0.40236040279000002"

$ws.Range("C46").Value = [double]"97.041797590000002"
$ws.Range("D46").Value = [double]"97.460478365"
$ws.Range("E46").Value = [double]"-0.4186807743"

$ws.Range("C47").Value = [double]"98.501512758000004"
$ws.Range("D47").Value = [double]"96.940833552000001"
$ws.Range("E47").Value = [double]"1.5606792053"

$ws.Range("C48").Value = [double]"96.775541509999996"
$ws.Range("D48").Value = [double]"96.983636528999995"
$ws.Range("E48").Value = [double]"-0.20809501960999999"

$ws.Range("C49").Value = [double]"97.696307489000006"
$ws.Range("D49").Value = [double]"97.723528970999993"
$ws.Range("E49").Value = [double]"-2.7221482158999999E-2"

$ws.Range("C50").Value = [double]"98.547778953999995"
$ws.Range("D50").Value = [double]"98.995585095999999"
$ws.Range("E50").Value = [double]"-0.44780614206000002"

$ws.Range("C51").Value = [double]"99.073255919999994"
$ws.Range("D51").Value = [double]"98.621587374000001"
$ws.Range("E51").Value = [double]"0.45166854650999999"

$ws.Range("C52").Value = [double]"98.693963736000001"
$ws.Range("D52").Value = [double]"98.578322159999999"
$ws.Range("E52").Value = [double]"0.11564157606"

$ws.Range("C53").Value = [double]"99.829087755000003"
$ws.Range("D53").Value = [double]"99.180575309999995"
$ws.Range("E53").Value = [double]"0.64851244425999999"

$ws.Range("C54").Value = [double]"100.09266555000001"
$ws.Range("D54").Value = [double]"100.31969019"
$ws.Range("E54").Value = [double]"-0.22702463140000001"

$ws.Range("C55").Value = [double]"100.40464325000001"
$ws.Range("D55").Value = [double]"100.07765381"
$ws.Range("E55").Value = [double]"0.32698944314"
